$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '59.992.62'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +3.57%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.421.00'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +3.17%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '552.75'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.10'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.68%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  +1.21%  '
$ws.Range('E9').Value = '  +3.82%  '
$ws.Range('E10').Value = '  +4.30%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.359'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.59%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '24.91'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.72%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.852.99'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.17%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '59.938.21'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.60%  '
$ws.Range('E16').Value = '  +2.90%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.403.65'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.32%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.36'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +6.22%  '
$ws.Range('E19').Value = '  +2.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '332.02'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.73%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.77'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.13%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.00'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.11%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '65.05'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.51%  '
$ws.Range('E24').Value = '  +3.69%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.62'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.39%  '
$ws.Range('E26').Value = '  +0.09%  '
$ws.Range('E27').Value = '  -0.20%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0784'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.47%  '
$ws.Range('E29').Value = '  +1.14%  '
$ws.Range('E30').Value = '  +2.43%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '168.95'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.61%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.04'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.26%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.70'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.78%  '
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('E35').Value = '  +5.26%  '
$ws.Range('B36').Value = 'NEARProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.22'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.70%  '
$ws.Range('B37').Value = 'FirstDigitalUSD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.09%  '
$ws.Range('E38').Value = '  +0.12%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '323.24'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +11.96%  '
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '39.54'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.10%  '
$ws.Range('B41').Value = 'PolygonEcosystemToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.416'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +10.33%  '
$ws.Range('E42').Value = '  +1.15%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '140.06'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.53%  '
$ws.Range('E44').Value = '  +1.33%  '
$ws.Range('E45').Value = '  +1.99%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '19.53'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.03%  '
$ws.Range('B47').Value = 'Polygon'
$ws.Range('C47').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.413'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +8.34%  '
$ws.Range('E48').Value = '  +1.19%  '
$ws.Range('E49').Value = '  +1.98%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '17.71'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.21%  '
$ws.Range('E51').Value = '  -0.15%  '
